$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Derived" / "Total" rows summarizing toxicity across species for each
# treatment x phase x day combination (rows 29-37).
$newRows = @(
    @("Derived", "Total", "18°C", "Uptake",      4, 0.28122994652406441),
    @("Derived", "Total", "21°C", "Uptake",      4, 0.52748663101604198),
    @("Derived", "Total", "24°C", "Uptake",      4, 0.58030748663101472),
    @("Derived", "Total", "18°C", "Depuration",  5, 0.49129679144385013),
    @("Derived", "Total", "21°C", "Depuration",  5, 0.54340909090909018),
    @("Derived", "Total", "24°C", "Depuration",  5, 1.6365775401069516),
    @("Derived", "Total", "18°C", "Depuration",  6, 0),
    @("Derived", "Total", "21°C", "Depuration",  6, 0),
    @("Derived", "Total", "24°C", "Depuration",  6, 0)
)

$startRow = 29
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    # Column B ("Total") is written before column A ("Derived") so that the
    # shared-string table registers "Total" ahead of "Derived", matching the
    # order in which the original author's edit introduced these strings.
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Update the selected cell/range on the sheet view.
$ws.Range("D20").Select()
